$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 839.7143
$ws.Range("I32").Value = 599.5
$ws.Range("J32").Value = 935.8
$ws.Range("K32").Value = 599.5
$ws.Range("L32").Value = 935.8
$ws.Range("M32").Value = -273.5
$ws.Range("N32").Value = -1587.8
$ws.Range("H39").Value = 2464.923
$ws.Range("I39").Value = 715.8570999999999
$ws.Range("K39").Value = 2147.5713
$ws.Range("M39").Value = -1851.5713
$ws.Range("H43").Value = 6800
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H97").Value = 5092.6665
$ws.Range("I97").Value = 3889
$ws.Range("J97").Value = 5333.4
$ws.Range("K97").Value = 11667
$ws.Range("L97").Value = 16000.2
$ws.Range("M97").Value = -11171
$ws.Range("N97").Value = -16992.2
$ws.Range("H100").Value = 2247.0952
$ws.Range("I100").Value = 2027.5
$ws.Range("K100").Value = 2027.5
$ws.Range("M100").Value = -1486.5
$ws.Range("H101").Value = 7391
$ws.Range("I101").Value = 7850
$ws.Range("J101").Value = 7276.25
$ws.Range("K101").Value = 23550
$ws.Range("L101").Value = 21828.75
$ws.Range("M101").Value = -21928
$ws.Range("N101").Value = -25072.75
$ws.Range("H107").Value = 1065.3334
$ws.Range("I107").Value = 1065.3334
$ws.Range("K107").Value = 1065.3334
$ws.Range("M107").Value = 854.6666
$ws.Range("H138").Value = 3214.5833
$ws.Range("I138").Value = 2735
$ws.Range("J138").Value = 8490
$ws.Range("K138").Value = 8205
$ws.Range("L138").Value = 25470
$ws.Range("M138").Value = -3065
$ws.Range("N138").Value = -35750

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3843.2683
$ws.Range("I32").Value = 2688.7896
$ws.Range("K32").Value = 2688.7896
$ws.Range("M32").Value = -2401.7896
$ws.Range("H45").Value = 5832.4736
$ws.Range("I45").Value = 6908.154
$ws.Range("J45").Value = 3501.8333
$ws.Range("K45").Value = 6908.154
$ws.Range("L45").Value = 3501.8333
$ws.Range("M45").Value = -6531.154
$ws.Range("N45").Value = -4255.8333
$ws.Range("H61").Value = 35717196
$ws.Range("I61").Value = 52633080
$ws.Range("K61").Value = 52633080
$ws.Range("M61").Value = -52632868
$ws.Range("H136").Value = 35717196
$ws.Range("I136").Value = 52633080
$ws.Range("K136").Value = 157899240
$ws.Range("M136").Value = -157896690

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1163.2
$ws.Range("I20").Value = 918.9
$ws.Range("J20").Value = 1651.8
$ws.Range("K20").Value = 918.9
$ws.Range("L20").Value = 1651.8
$ws.Range("M20").Value = -671.9
$ws.Range("N20").Value = -2145.8
$ws.Range("H99").Value = 2041.826
$ws.Range("I99").Value = 1907
$ws.Range("J99").Value = 2188.9092
$ws.Range("K99").Value = 1907
$ws.Range("L99").Value = 2188.9092
$ws.Range("M99").Value = -409
$ws.Range("N99").Value = -5184.9092

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4991.143
$ws.Range("I31").Value = 3450.9375
$ws.Range("J31").Value = 9919.799999999999
$ws.Range("K31").Value = 3450.9375
$ws.Range("L31").Value = 9919.799999999999
$ws.Range("M31").Value = -3155.9375
$ws.Range("N31").Value = -10509.8
$ws.Range("H34").Value = 4991.143
$ws.Range("I34").Value = 3450.9375
$ws.Range("J34").Value = 9919.799999999999
$ws.Range("K34").Value = 3450.9375
$ws.Range("L34").Value = 9919.799999999999
$ws.Range("M34").Value = -3248.9375
$ws.Range("N34").Value = -10323.8

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 756.75
$ws.Range("I63").Value = 342.33334
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1027.00002
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -278.0000199999999
$ws.Range("N63").Value = -7498
$ws.Range("H66").Value = 756.75
$ws.Range("I66").Value = 342.33334
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 3081.00006
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = 662.9999399999997
$ws.Range("N66").Value = -25488
$ws.Range("H75").Value = 385.4
$ws.Range("I75").Value = 400
$ws.Range("K75").Value = 1200
$ws.Range("M75").Value = -202
$ws.Range("H78").Value = 385.4
$ws.Range("I78").Value = 400
$ws.Range("K78").Value = 3600
$ws.Range("M78").Value = 1392
$ws.Range("H81").Value = 306725
$ws.Range("I81").Value = 470000
$ws.Range("K81").Value = 1410000
$ws.Range("M81").Value = -1408877
$ws.Range("H84").Value = 306725
$ws.Range("I84").Value = 470000
$ws.Range("K84").Value = 4230000
$ws.Range("M84").Value = -4224384
$ws.Range("H103").Value = 1348.8182
$ws.Range("I103").Value = 283.42856
$ws.Range("J103").Value = 3213.25
$ws.Range("K103").Value = 850.28568
$ws.Range("L103").Value = 9639.75
$ws.Range("M103").Value = 28.71432000000004
$ws.Range("N103").Value = -11397.75
$ws.Range("H109").Value = 1424.6
$ws.Range("I109").Value = 1424.6
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 4273.799999999999
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -3233.799999999999
$ws.Range("N109").ClearContents()
$ws.Range("H117").Value = 1682.9
$ws.Range("I117").Value = 351
$ws.Range("J117").Value = 2253.7144
$ws.Range("K117").Value = 1053
$ws.Range("L117").Value = 6761.1432
$ws.Range("M117").Value = 2389
$ws.Range("N117").Value = -13645.1432
$ws.Range("H123").Value = 899
$ws.Range("J123").Value = 899
$ws.Range("L123").Value = 2697
$ws.Range("N123").Value = -7597

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 151.5
$ws.Range("J2").Value = 134.4
$ws.Range("L2").Value = 134.4
$ws.Range("N2").Value = -360.4
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H46").Value = 7298
$ws.Range("I46").Value = 1625
$ws.Range("J46").Value = 29990
$ws.Range("K46").Value = 1625
$ws.Range("L46").Value = 29990
$ws.Range("M46").Value = -1469
$ws.Range("N46").Value = -30302
$ws.Range("H70").Value = 5830.2
$ws.Range("I70").Value = 6801.3335
$ws.Range("J70").Value = 5414
$ws.Range("K70").Value = 6801.3335
$ws.Range("L70").Value = 5414
$ws.Range("M70").Value = -6531.3335
$ws.Range("N70").Value = -5954
$ws.Range("H73").Value = 5830.2
$ws.Range("I73").Value = 6801.3335
$ws.Range("J73").Value = 5414
$ws.Range("K73").Value = 6801.3335
$ws.Range("L73").Value = 5414
$ws.Range("M73").Value = -5865.3335
$ws.Range("N73").Value = -7286
$ws.Range("H80").Value = 1991
$ws.Range("I80").Value = 1702.7142
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 1702.7142
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -704.7141999999999
$ws.Range("N80").Value = -4996
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 1991
$ws.Range("I83").Value = 1702.7142
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 8513.571
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -3521.571
$ws.Range("N83").Value = -24984
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 600.11536
$ws.Range("I55").Value = 468
$ws.Range("K55").Value = 468
$ws.Range("M55").Value = -295
$ws.Range("H82").Value = 2150.8948
$ws.Range("I82").Value = 2059.8
$ws.Range("J82").Value = 2252.111
$ws.Range("K82").Value = 2059.8
$ws.Range("L82").Value = 2252.111
$ws.Range("M82").Value = -1698.8
$ws.Range("N82").Value = -2974.111
$ws.Range("H85").Value = 2150.8948
$ws.Range("I85").Value = 2059.8
$ws.Range("J85").Value = 2252.111
$ws.Range("K85").Value = 2059.8
$ws.Range("L85").Value = 2252.111
$ws.Range("M85").Value = -811.8000000000002
$ws.Range("N85").Value = -4748.111

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 45486.086
$ws.Range("I81").Value = 60108.176
$ws.Range("J81").Value = 4056.8333
$ws.Range("K81").Value = 120216.352
$ws.Range("L81").Value = 8113.6666
$ws.Range("M81").Value = -119155.352
$ws.Range("N81").Value = -10235.6666
$ws.Range("H84").Value = 45486.086
$ws.Range("I84").Value = 60108.176
$ws.Range("J84").Value = 4056.8333
$ws.Range("K84").Value = 601081.76
$ws.Range("L84").Value = 40568.333
$ws.Range("M84").Value = -595777.76
$ws.Range("N84").Value = -51176.333
$ws.Range("H96").Value = 3455.04
$ws.Range("I96").Value = 1484.1428
$ws.Range("K96").Value = 1484.1428
$ws.Range("M96").Value = -111.1428000000001
$ws.Range("H100").Value = 2932.5715
$ws.Range("I100").Value = 2905.6
$ws.Range("K100").Value = 5811.2
$ws.Range("M100").Value = -5270.2
$ws.Range("H122").Value = 1747.1818
$ws.Range("I122").Value = 1747.1818
$ws.Range("K122").Value = 5241.5454
$ws.Range("M122").Value = -2791.5454
